$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the repeated values in A3:A4 and B3:B4 (they'll be represented by the merge)
$ws.Range("A3:A4").ClearContents()
$ws.Range("B3:B4").ClearContents()

# Merge the cells
$ws.Range("A2:A4").MergeCells = $true
$ws.Range("B2:B4").MergeCells = $true

# Update the selection to match target
$ws.Range("F13").Select()
